$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# Update the ValueSet URL (row 2) from pythia -> cicada
$meta.Range("B2").Value = "http://fhirfli.dev/fhir/ig/cicada/ValueSet/eval-reason"

# Update the Date value (row 8)
$meta.Range("B8").Value = "2026-02-11T14:37:07-05:00"

# Insert a new "Jurisdiction" row right after "Contact" (row 10), pushing
# Description/Purpose/Copyright/Immutable down by one row.
$meta.Rows.Item(11).Insert()
$meta.Range("A11").Value = "Jurisdiction"
$meta.Range("B11").Value = ""
$meta.Range("A11:B11").Style = $meta.Range("A12:B12").Style

# --- Include sheet ----------------------------------------------------
$inc = $wb.Worksheets.Item("Include from Evaluation Reaso")

# Update the CodeSystem URI (last row) from pythia -> cicada
$inc.Range("B18").Value = "http://fhirfli.dev/fhir/ig/cicada/CodeSystem/EvalReason"

# Rename the sheet
$inc.Name = "Include #0"
